$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before K (shifts old K "Prioridad" -> L, old L "Fecha solución" -> M)
$ws.Columns("K").Insert()

# Rename the "Descripción" header (column C) to "Reporte"
$ws.Range("C1").Value = "Reporte"

# Give the newly inserted column a header of "Severidad"
$ws.Range("K1").Value = "Severidad"

# Narrow column B a bit (closest achievable approximation of the author's manual resize)
$ws.Columns("B").ColumnWidth = 26.86

# Leave the selection on the newly edited cell, as in the saved file
[void]$ws.Range("K2").Select()
